$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the September row label (through 09-05 -> through 09-06)
$ws.Range("A10").Value = "September (through 09-06)"

# Update September row values (row 10)
$ws.Range("B10").Value = 5
$ws.Range("C10").Value = 11
$ws.Range("D10").Value = 14
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = 13
$ws.Range("G10").Value = 21
$ws.Range("H10").Value = 26

# Update Total row values (row 11)
$ws.Range("B11").Value = 199
$ws.Range("C11").Value = 392
$ws.Range("D11").Value = 565
$ws.Range("E11").Value = 499
$ws.Range("F11").Value = 362
$ws.Range("G11").Value = 805
$ws.Range("H11").Value = 1097
